# MQuiz4 and slides for PL
# Shifts the "topic"/"reading" schedule content in rows 18-23 of the
# Calendar2021 sheet down by one slot to make room for a new
# "Reinforcement Learning 3" topic row, inserts the updated Prop Logic
# slide reference, adds the Quiz4/Quiz3-retake note, and moves the
# active-cell selection from D20 to A20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calendar2021")

# Row 18: was "Proposition Logic and Logical Agents" / "Chp 7.1-7.4"
# becomes the new "Reinforcement Learning 3" topic (reading cleared).
$ws.Range("D18").Value = "Reinforcement Learning 3;video(https://canvas.jmu.edu/courses/1775272/modules)"
$ws.Range("E18").Value = ""

# Row 19: was "Theorm Proving and Resolution" / "Chp 7.5"
# becomes the updated Proposition Logic topic (with slides added),
# the old "Chp 7.1-7.4" reading, and a new Quiz4/Quiz3-retake note.
$ws.Range("D19").Value = "Proposition Logic and Logical Agentsslides(slides/13_PropLogic.pdf); video(https://canvas.jmu.edu/courses/1775272/modules)"
$ws.Range("E19").Value = "Chp 7.1-7.4"
$ws.Range("F19").Value = "Quiz 4(mquizzes/mquiz4/mquiz4.php);Quiz3 Retake"

# Row 20: was "FOL, Unification"
# becomes "Theorm Proving and Resolution" / "Chp 7.5"
$ws.Range("D20").Value = "Theorm Proving and Resolution"
$ws.Range("E20").Value = "Chp 7.5"

# Row 21: was "FOL Resolution and Chaining"
# becomes "FOL, Unification"
$ws.Range("D21").Value = "FOL, Unification"

# Row 22: was "Prob"
# becomes "FOL Resolution and Chaining"
$ws.Range("D22").Value = "FOL Resolution and Chaining"

# Row 23: was empty
# becomes "Prob"
$ws.Range("D23").Value = "Prob"

# Move the saved selection/active cell from D20 to A20.
$ws.Range("A20").Select()
